$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (Name, Position, Team) reflecting the final values
$rows = @{
    3  = @("Devin Vassell", "SG,SF", "San Antonio Spurs")
    6  = @("Desmond Bane", "SG,SF", "Memphis Grizzlies")
    7  = @("Patrick Williams", "PF", "Chicago Bulls")
    8  = @("Norman Powell", "SG,SF", "LA Clippers")
    9  = @("LeBron James", "SF,PF", "Los Angeles Lakers")
    10 = @("Goga Bitadze", "C", "Orlando Magic")
    11 = @("Walker Kessler", "C", "Utah Jazz")
    13 = @("Brandin Podziemski", "SG", "Golden State Warriors")
    14 = @("Devin Booker", "PG,SG", "Phoenix Suns")
    15 = @("Alperen Sengün", "C", "Houston Rockets")
    16 = @("Shaedon Sharpe", "SG,SF", "Portland Trail Blazers")
}

foreach ($r in $rows.Keys) {
    $entry = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
}
